$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'compression underwear women'
$ws.Range("A2").Value = 'compression underwear women high waist'
$ws.Range("A3").Value = 'compression underwear women workout'
$ws.Range("A4").Value = 'compression upper leg'
$ws.Range("A5").Value = 'compression waist'
$ws.Range("A6").Value = 'compression waist women'
$ws.Range("A7").Value = 'compression waist wrap'
$ws.Range("A8").Value = 'compression wear'
$ws.Range("A9").Value = 'compression wear for women'
$ws.Range("A10").Value = 'compression wears'
$ws.Range("A11").Value = 'compression winter'
$ws.Range("A12").Value = 'compression woman'
$ws.Range("A13").Value = 'compression woman leggings'
$ws.Range("A14").Value = 'compression women'
$ws.Range("A15").Value = 'compression women calf'
$ws.Range("A16").Value = 'compression women capri'
$ws.Range("A17").Value = 'compression women circulation'
$ws.Range("A18").Value = 'compression women cwx'
$ws.Range("A19").Value = 'compression women knee'
$ws.Range("A20").Value = 'compression women leggings'
$ws.Range("A21").Value = 'compression women pants'
$ws.Range("A22").Value = 'compression women running'
$ws.Range("A23").Value = 'compression women shorts'
$ws.Range("A24").Value = 'compression women tights'
$ws.Range("A25").Value = 'compression women top'
$ws.Range("A26").Value = 'compression women tops'
$ws.Range("A27").Value = 'compression women underwear'
$ws.Range("A28").Value = 'compression womens'
$ws.Range("A29").Value = 'compression workout'
$ws.Range("A30").Value = 'compression workout capri'
$ws.Range("A31").Value = 'compression workout capris for women'
$ws.Range("A32").Value = 'compression workout clothes'
$ws.Range("A33").Value = 'compression workout leggings women'
$ws.Range("A34").Value = 'compression workout pants for women'
$ws.Range("A35").Value = 'compression workout pants women'
$ws.Range("A36").Value = 'compression workout tights women'
$ws.Range("A37").Value = 'compression workout tops women'
$ws.Range("A38").Value = 'compression wrap waist'
$ws.Range("A39").Value = 'compression wrap women'
$ws.Range("A40").Value = 'compression x'
$ws.Range("A41").Value = 'compression yoga'
$ws.Range("A42").Value = 'compression yoga capris women'
$ws.Range("A43").Value = 'compression yoga leggings'
$ws.Range("A44").Value = 'compression yoga pant'
$ws.Range("A45").Value = 'compression yoga pants'
$ws.Range("A46").Value = 'compression yoga pants for women'
$ws.Range("A47").Value = 'compression yoga pants women'
$ws.Range("A48").Value = 'compression yoga pants women high waist'
$ws.Range("A49").Value = 'compression yoga tights'
$ws.Range("A50").Value = 'compression z leggings'
$ws.Range("A51").Value = 'compressions knee'
$ws.Range("A52").Value = 'compressions shorts'
$ws.Range("A53").Value = 'compressions tights for men'
$ws.Range("A54").Value = 'compressions underwear for women'
$ws.Range("A55").Value = 'compressionz leggings women'
$ws.Range("A56").Value = 'compressionz women''s compression pants'
$ws.Range("A57").Value = 'compresson shorts'
$ws.Range("A58").Value = 'compresson shorts men'
$ws.Range("A59").Value = 'compretion tights'
$ws.Range("A60").Value = 'comression leggings'
$ws.Range("A61").Value = 'concrete compression testing machine'
$ws.Range("A62").Value = 'conditioning equipment'
$ws.Range("A63").Value = 'confort women'
$ws.Range("A64").Value = 'conpression garments'
$ws.Range("A65").Value = 'construction back brace'
$ws.Range("A66").Value = 'construction clothes for women'
$ws.Range("A67").Value = 'contour cool toned'
$ws.Range("A68").Value = 'contour thermal'
$ws.Range("A69").Value = 'contractions game'
$ws.Range("A70").Value = 'control tape'
$ws.Range("A71").Value = 'control tights'
$ws.Range("A72").Value = 'control top black tights'
$ws.Range("A73").Value = 'control top high waist tights'
$ws.Range("A74").Value = 'control top leggings black'
$ws.Range("A75").Value = 'control top red leggings'
$ws.Range("A76").Value = 'control top yoga pants'
$ws.Range("A77").Value = 'cool compression tights'
$ws.Range("A78").Value = 'cool gear compression'
$ws.Range("A79").Value = 'cool gifts for runners'
$ws.Range("A80").Value = 'cool leggings for women'
$ws.Range("A81").Value = 'cool pants for women'
$ws.Range("A82").Value = 'cool pants women'
$ws.Range("A83").Value = 'cool ski pants'
$ws.Range("A84").Value = 'cool tights for women'
$ws.Range("A85").Value = 'cool tights women'
$ws.Range("A86").Value = 'cool weather running pants'
$ws.Range("A87").Value = 'coolmax leggings'
$ws.Range("A88").Value = 'coolmax shirts'
$ws.Range("A89").Value = 'coolmax shirts women'
$ws.Range("A90").Value = 'coolmax underwear men'
$ws.Range("A91").Value = 'copper back support for women'
$ws.Range("A92").Value = 'copper clothing women'
$ws.Range("A93").Value = 'copper compression'
$ws.Range("A94").Value = 'copper compression back'
$ws.Range("A95").Value = 'copper compression for knee'
$ws.Range("A96").Value = 'copper compression for knees'
$ws.Range("A97").Value = 'copper compression for women'
$ws.Range("A98").Value = 'copper compression leggings women'
$ws.Range("A99").Value = 'copper compression leggings women'
$ws.Range("A100").Value = 'copper compression pants'
